# Add experiment C sample dates (column I, "sample_date") for rows 34-65,
# plus default 0 values for number_arms_dropped/arm_twist/arms_crossed
# (columns O/P/Q) on rows that did not already carry those counts.
# Also appends a new (otherwise empty) row 66 that only carries the date
# style on I66, and updates the active selection to K58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> sample_date (Excel serial date number), matching column H's
# "inoc_date" style (s="1") exactly by copying H's formatting into I.
$sampleDates = [ordered]@{
  34 = 44733; 35 = 44733; 36 = 44733; 37 = 44733; 38 = 44733; 39 = 44733;
  40 = 44733; 41 = 44733; 42 = 44733; 43 = 44733; 44 = 44733; 45 = 44733;
  46 = 44741; 47 = 44741; 48 = 44742; 49 = 44743; 50 = 44744; 51 = 44744;
  52 = 44744; 53 = 44744; 54 = 44744; 55 = 44745; 56 = 44745; 57 = 44745;
  58 = 44745; 59 = 44746; 60 = 44746; 61 = 44747; 62 = 44747; 63 = 44748;
  64 = 44748; 65 = 44748
}

# Rows that previously had no number_arms_dropped/twist/crossed values at
# all and now get explicit zeros in O/P/Q.
$zeroArmCols = @(34,35,36,37,38,39,40,41,42,43,44,45,50,51,56,57,59,61,63)

foreach ($row in $sampleDates.Keys) {
    $iCell = $ws.Range("I$row")
    # Copy H's cell formatting (date style) into I before writing the value,
    # so I ends up with the same number format / style index as H.
    $ws.Range("H$row").Copy($iCell)
    $iCell.Value = $sampleDates[$row]

    if ($zeroArmCols -contains $row) {
        $ws.Range("O$row").Value = 0
        $ws.Range("P$row").Value = 0
        $ws.Range("Q$row").Value = 0
    }
}

# New trailing row 66: a single styled-but-empty date cell at I66.
$ws.Range("H34").Copy($ws.Range("I66"))
$ws.Range("I66").ClearContents()

# Match the author's final selection/viewport.
$ws.Range("K58").Select()
